# Update "想去人数" (interest count) figures in the "展览", "演出" and
# "全部类型" sheets to the latest scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 301
$ws1.Range("F3").Value = 13836
$ws1.Range("F4").Value = 258
$ws1.Range("F6").Value = 181
$ws1.Range("F7").Value = 276
$ws1.Range("F8").Value = 494
$ws1.Range("F14").Value = 450
$ws1.Range("F15").Value = 5796
$ws1.Range("F16").Value = 128
$ws1.Range("F17").Value = 87
$ws1.Range("F18").Value = 977
$ws1.Range("F19").Value = 95
$ws1.Range("F20").Value = 52
$ws1.Range("F21").Value = 150
$ws1.Range("F22").Value = 241

# --- Sheet: 演出 ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14

# --- Sheet: 全部类型 -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 301
$ws4.Range("F3").Value = 13836
$ws4.Range("F4").Value = 258
$ws4.Range("F6").Value = 181
$ws4.Range("F7").Value = 276
$ws4.Range("F8").Value = 494
$ws4.Range("F14").Value = 450
$ws4.Range("F15").Value = 5796
$ws4.Range("F16").Value = 128
$ws4.Range("F17").Value = 87
$ws4.Range("F18").Value = 977
$ws4.Range("F19").Value = 95
$ws4.Range("F20").Value = 52
$ws4.Range("F21").Value = 150
$ws4.Range("F22").Value = 241
$ws4.Range("F23").Value = 14
